$d = $word.ActiveDocument

$pairs = @(
    @("57×89=", "70×66="),
    @("48×25=", "54×44="),
    @("96×61=", "58×59="),
    @("32×68=", "18×86="),
    @("40×11=", "75×88="),
    @("20×47=", "25×38="),
    @("21×66=", "56×59="),
    @("83×70=", "60×74="),
    @("85×99=", "76×21="),
    @("97×66=", "66×44="),
    @("51×78=", "88×48="),
    @("30×12=", "63×59="),
    @("36×48=", "90×82="),
    @("56×80=", "15×86="),
    @("31×19=", "55×13="),
    @("44×57=", "18×21="),
    @("81×84=", "11×28="),
    @("59×86=", "90×90="),
    @("48×56=", "98×64="),
    @("47×62=", "66×74="),
    @("80×47=", "15×86="),
    @("98×81=", "42×44="),
    @("23×20=", "22×90="),
    @("57×27=", "99×21="),
    @("87×73=", "72×46=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
